$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$e2 = @'

                    <filter>
                        <interfaces xmlns="http://openconfig.net/yang/interfaces">
                            <interface>
                                <name>GigabitEthernet0/3/2</name>
                                <subinterfaces>
                                    <subinterface>
                                        <index>100</index>
                                    </subinterface>
                                </subinterfaces>
                            </interface>
                        </interfaces>
                    </filter>
                     
-------------------
                    <filter>
                        <network-instances xmlns="http://openconfig.net/yang/network-instance">
                        <network-instance>
                        <name>Prueba_LxVPN</name>
                        </network-instance>
                        </network-instances>
                    </filter>
                    
'@
$f2 = @'
<rpc-reply message-id="urn:uuid:3eef521f-7de7-4a12-a7d8-fe82f1370032">
  <data>
    <interfaces>
      <interface>
        <name>GigabitEthernet0/3/2</name>
      </interface>
    </interfaces>
  </data>
</rpc-reply>
 
-------------------
<rpc-reply message-id="urn:uuid:2fda1e91-e66f-484d-92c2-dfb610339a5b">
  <data>
    <network-instances>
      <network-instance>
        <name>Prueba_LxVPN</name>
        <config>
          <name>Prueba_LxVPN</name>
          <type>oc-ni-types:L3VRF</type>
        </config>
        <interfaces>
          <interface>
            <id>GigabitEthernet0/3/2</id>
            <config>
              <id>GigabitEthernet0/3/2</id>
              <interface>GigabitEthernet0/3/2</interface>
              <subinterface>0</subinterface>
            </config>
          </interface>
        </interfaces>
        <protocols>
          <protocol>
            <identifier>oc-pol-types:STATIC</identifier>
            <name>default</name>
            <config>
              <identifier>oc-pol-types:STATIC</identifier>
              <name>default</name>
            </config>
          </protocol>
          <protocol>
            <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>
            <name>default</name>
            <config>
              <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>
              <name>default</name>
            </config>
          </protocol>
        </protocols>
      </network-instance>
    </network-instances>
  </data>
</rpc-reply>

'@
$g2 = @'

  <edit-config>
    <target>
      <candidate/>
    </target>
    <config>
      <interfaces xmlns="http://openconfig.net/yang/interfaces">
        <interface>
          <name>GigabitEthernet0/3/2</name>
          <subinterfaces>
            <subinterface>
              <index>100</index>
              <config>
                <index>100</index>
              </config>
            </subinterface>
          </subinterfaces>
        </interface>
      </interfaces>
    </config>
  </edit-config> 
-------------------
  <edit-config>
    <target>
      <candidate/>
    </target>
    <config>
      <network-instances xmlns="http://openconfig.net/yang/network-instance">
        <network-instance>
          <name>Prueba_LxVPN</name>
          <config>
            <name>Prueba_LxVPN</name>
            <type xmlns:oc-ni-types="http://openconfig.net/yang/network-instance-types">oc-ni-types:L3VRF</type>
          </config>
          <interfaces>
            <interface>
              <id>GigabitEthernet0/3/2</id>
              <config>
              <id>GigabitEthernet0/3/2</id>
                <interface>GigabitEthernet0/3/2</interface>
                <subinterface>100</subinterface>
              </config>
            </interface>
          </interfaces>
        </network-instance>
      </network-instances>
    </config>
  </edit-config>
'@
$h2 = @'
- Response of edit-config: <rpc-reply message-id="urn:uuid:800cdfdc-b94d-49be-a13f-2d3fe4044913">
  <ok/>
</rpc-reply>
 

 - Response of commit: <rpc-reply xmlns:nc-ext="urn:huawei:yang:huawei-ietf-netconf-ext" message-id="urn:uuid:20f57bc4-1ad6-44eb-9b9e-be8804dbce4c" nc-ext:flow-id="252">
  <ok/>
</rpc-reply>
 
-------------------

'@
$i2 = @'
<rpc-reply message-id="urn:uuid:dd8d79fd-8b63-4b44-bd2e-7fcf9c3a3dfd">
  <data>
    <interfaces>
      <interface>
        <name>GigabitEthernet0/3/2</name>
        <subinterfaces>
          <subinterface>
            <index>100</index>
            <config>
              <index>100</index>
              <enabled>true</enabled>
            </config>
            <ipv4>
              <config>
                <enabled>true</enabled>
              </config>
            </ipv4>
          </subinterface>
        </subinterfaces>
      </interface>
    </interfaces>
  </data>
</rpc-reply>
 
-------------------

'@

$ws.Range("E2").Value2 = $e2
$ws.Range("F2").Value2 = $f2
$ws.Range("G2").Value2 = $g2
$ws.Range("H2").Value2 = $h2
$ws.Range("I2").Value2 = $i2
